$wb = $excel.ActiveWorkbook

# ALC!row34
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3762.6667
$ws.Range("I34").Value = 3762.6667
$ws.Range("K34").Value = 3762.6667
$ws.Range("M34").Value = -3559.6667

# ALC!row36
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 3762.6667
$ws.Range("I36").Value = 3762.6667
$ws.Range("K36").Value = 3762.6667
$ws.Range("M36").Value = -3047.6667

# ALC!row47
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 15000
$ws.Range("I47").Value = 10000
$ws.Range("J47").Value = 20000
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 20000
$ws.Range("M47").Value = -9028
$ws.Range("N47").Value = -21944

# ALC!row54
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 7966.6665
$ws.Range("I54").Value = 7966.6665
$ws.Range("K54").Value = 7966.6665
$ws.Range("M54").Value = -7480.6665

# ALC!row62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3721.6667
$ws.Range("I62").Value = 1652.0834
$ws.Range("K62").Value = 1652.0834
$ws.Range("M62").Value = -1028.0834

# ALC!row65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3721.6667
$ws.Range("I65").Value = 1652.0834
$ws.Range("K65").Value = 8260.416999999999
$ws.Range("M65").Value = -5140.416999999999

# ALC!row80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 657.1818
$ws.Range("I80").Value = 719.4
$ws.Range("K80").Value = 2158.2
$ws.Range("M80").Value = -1160.2

# ALC!row83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 657.1818
$ws.Range("I83").Value = 719.4
$ws.Range("K83").Value = 6474.599999999999
$ws.Range("M83").Value = -1482.599999999999

# ALC!row98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 668.9286
$ws.Range("I98").Value = 505.81818
$ws.Range("K98").Value = 505.81818
$ws.Range("M98").Value = 992.18182

# ALC!row122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 668.9286
$ws.Range("I122").Value = 505.81818
$ws.Range("K122").Value = 1517.45454
$ws.Range("M122").Value = 932.54546

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2943670.5
$ws.Range("I32").Value = 563.75
$ws.Range("K32").Value = 563.75
$ws.Range("M32").Value = -276.75

# ARM!row43
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 10037500
$ws.Range("I43").Value = 10037500
$ws.Range("K43").Value = 10037500
$ws.Range("M43").Value = -10037187

# ARM!row76
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20676

# ARM!row79
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22340

# ARM!row88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 829.8182
$ws.Range("I88").Value = 379.5
$ws.Range("J88").Value = 1370.2
$ws.Range("K88").Value = 379.5
$ws.Range("L88").Value = 1370.2
$ws.Range("M88").Value = 26.5
$ws.Range("N88").Value = -2182.2

# ARM!row91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 829.8182
$ws.Range("I91").Value = 379.5
$ws.Range("J91").Value = 1370.2
$ws.Range("K91").Value = 379.5
$ws.Range("L91").Value = 1370.2
$ws.Range("M91").Value = 1024.5
$ws.Range("N91").Value = -4178.2

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2085.2
$ws.Range("I97").Value = 2430
$ws.Range("J97").Value = 1740.4
$ws.Range("K97").Value = 2430
$ws.Range("L97").Value = 1740.4
$ws.Range("M97").Value = -1934
$ws.Range("N97").Value = -2732.4

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 7356431
$ws.Range("I102").Value = 11364484
$ws.Range("K102").Value = 11364484
$ws.Range("M102").Value = -11362862

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3090.1667
$ws.Range("I122").Value = 2166.6667
$ws.Range("K122").Value = 6500.000100000001
$ws.Range("M122").Value = -4050.000100000001

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2584.476
$ws.Range("I86").Value = 1269.6154
$ws.Range("J86").Value = 4721.125
$ws.Range("K86").Value = 1269.6154
$ws.Range("L86").Value = 4721.125
$ws.Range("M86").Value = -146.6153999999999
$ws.Range("N86").Value = -6967.125

# BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2584.476
$ws.Range("I89").Value = 1269.6154
$ws.Range("J89").Value = 4721.125
$ws.Range("K89").Value = 6348.076999999999
$ws.Range("L89").Value = 23605.625
$ws.Range("M89").Value = -732.0769999999993
$ws.Range("N89").Value = -34837.625

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 577.8
$ws.Range("I94").Value = 575
$ws.Range("K94").Value = 575
$ws.Range("M94").Value = -124

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 100000910
$ws.Range("I99").Value = 111112010
$ws.Range("J99").Value = 1090
$ws.Range("K99").Value = 111112010
$ws.Range("L99").Value = 1090
$ws.Range("M99").Value = -111110512
$ws.Range("N99").Value = -4086

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3923.12
$ws.Range("I134").Value = 1173.95
$ws.Range("J134").Value = 14919.8
$ws.Range("K134").Value = 3521.85
$ws.Range("L134").Value = 44759.39999999999
$ws.Range("M134").Value = -986.8500000000004
$ws.Range("N134").Value = -49829.39999999999

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3713.5151
$ws.Range("I31").Value = 2349.6924
$ws.Range("K31").Value = 2349.6924
$ws.Range("M31").Value = -2054.6924

# CRP!row33
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I33").Value = 2042.7142
$ws.Range("J33").Value = 8766.333000000001
$ws.Range("K33").Value = 2042.7142
$ws.Range("L33").Value = 8766.333000000001
$ws.Range("M33").Value = -1663.7142
$ws.Range("N33").Value = -9524.333000000001

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3713.5151
$ws.Range("I34").Value = 2349.6924
$ws.Range("K34").Value = 2349.6924
$ws.Range("M34").Value = -2147.6924

# CRP!row44
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 5000
$ws.Range("I44").Value = 5000
$ws.Range("K44").Value = 5000
$ws.Range("M44").Value = -4558

# CRP!row55
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 7718.6
$ws.Range("I55").Value = 4648.5
$ws.Range("J55").Value = 19999
$ws.Range("K55").Value = 4648.5
$ws.Range("L55").Value = 19999
$ws.Range("M55").Value = -4333.5
$ws.Range("N55").Value = -20629

# CUL!row32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2712.4285
$ws.Range("I131").Value = 1999.5
$ws.Range("K131").Value = 5998.5
$ws.Range("M131").Value = -958.5

# GSM!row44
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 4467.5
$ws.Range("I44").Value = 5000
$ws.Range("K44").Value = 5000
$ws.Range("M44").Value = -4404

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# GSM!row86
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# GSM!row89
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# GSM!row97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1162.8182
$ws.Range("J97").Value = 956.8333
$ws.Range("L97").Value = 956.8333
$ws.Range("N97").Value = -1948.8333

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2144.9
$ws.Range("I46").Value = 592
$ws.Range("J46").Value = 2981.077
$ws.Range("K46").Value = 592
$ws.Range("L46").Value = 2981.077
$ws.Range("M46").Value = -404
$ws.Range("N46").Value = -3357.077

# LTW!row100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4314.2144
$ws.Range("I100").Value = 4480
$ws.Range("K100").Value = 4480
$ws.Range("M100").Value = -3939

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3902
$ws.Range("I136").Value = 3145.0833
$ws.Range("J136").Value = 5415.8335
$ws.Range("K136").Value = 9435.249899999999
$ws.Range("L136").Value = 16247.5005
$ws.Range("M136").Value = -6885.249899999999
$ws.Range("N136").Value = -21347.5005
